# Daily attendance processing - 2025-10-08 10:47:12
# Swap the order of the two comma-separated entries in column G ("Recorded By")
# for every row where the value is exactly "System, dnasr281@gmail.com" or
# "admin@admin.com, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$startRow = $used.Row

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $startRow + $i
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Text

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "admin@admin.com, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, admin@admin.com"
    }
}
